# Rename the column header suffixes from "_old"/"_new" to the concrete
# format-version names "_FV2210"/"_FV2304", wrap the data range in an Excel
# Table (ListObject), and freeze the header row - matching the commit
# "adapt column header formatting to respective input file names".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:J hold the "<Label>_old" headers -> "<Label>_FV2210"
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
# Columns L:U hold the "<Label>_new" headers -> "<Label>_FV2304"
# (column K is the unchanged "diff" header in between)
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = $cell.Value2 -replace "_old$", "_FV2210"
}
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = $cell.Value2 -replace "_new$", "_FV2304"
}

# Turn the whole used range into a native Excel table (adds xl/tables/table1.xml,
# the worksheet <tableParts> entry and the sheet1 rels/content-type wiring).
$headerRange = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $headerRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"
# Source diff ships the table with no named style (plain banded rows only) -
# drop the auto-assigned default table style name.
$table.TableStyle = ""

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
